# Commit message: "Connected to the Database."
# The shared credentials sheet ("test") had its generated Password value
# (row 6, column C) refreshed to a newly generated value after the
# automation connected to the database and rotated the password.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")

# Previous value was "newPT_638*235"; it is replaced by the newly
# generated password "newPT_353*883".
$ws.Range("C6").Value = "newPT_353*883"
